$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("Season") before the existing table, shifting
# everything (including the Table1 structured table) one column to the right.
$ws.Columns("A:A").Insert()

# Resize/reposition the structured table to its new location.
$lo = $ws.ListObjects(1)
$lo.Resize($ws.Range("B1:N2"))

# New "Season" header + existing season's value.
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"

# Add a new table row for the 23/24 season and populate it.
$newRow = $lo.ListRows.Add()
$ws.Range("A3").Value = "23/24"
$ws.Range("B3").Value = "Nottingham"
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 49
$ws.Range("J3").Value = 67
$ws.Range("K3").Value = -18
$ws.Range("L3").Value = "Chris Wood"
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 78

# Apply a plain autofilter button to the new "Season" column (outside the table).
$null = $ws.Range("A1").AutoFilter()

# Record the (hidden) filter-database defined name Excel creates for the sheet
# whenever an autofilter is present.
$n = $ws.Names.Add("_xlnm._FilterDatabase", "='Nottingham Stats'!`$A`$1:`$A`$1")
$n.Visible = $false

# Update the selection to match the final edit position.
$null = $ws.Range("N3").Select()
